# Refresh the cryptocurrency price/volume snapshot (cols D and E) with
# the latest values pulled by the scheduled GitHub Actions scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.096.77"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.880.46"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5041"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3835"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08552"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.266"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "1.874.10"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.217"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.098"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "28.137.53"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.258"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").Value = "2.095.75"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1051"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.056"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.642"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.605"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.715"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02461"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06563"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6516"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("E41").Value = "  -7.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6178"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.301"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.025"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
